# Add season-record columns (Wins / Losses / Ties) to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the formatting of the last existing header cell
# (AC1) onto the three new header cells so they pick up the same bold /
# centered / bordered style (style index 1), then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-44: every player on the roster shares the team's season
# record, so the same Wins/Losses/Ties values repeat down the columns.
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 88
    $ws.Cells.Item($row, 31).Value = 74
    $ws.Cells.Item($row, 32).Value = 0
}
